$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$ws1 = $wb.Worksheets.Item(1)

# The original sheet had a duplicated "Contact" / "No display for ContactDetail" row
# (rows 10 and 11). Remove the duplicate row (row 11) so a new "Jurisdiction" row can
# take its place without growing the table, matching the new A1:B20 dimension.
$ws1.Rows.Item(11).Delete()

# Version bump
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Updated publication date
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank; now populated
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# Replace the (now de-duplicated) Contact row with a Jurisdiction row
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# --- Sheet 2: "Elements" ---
$ws2 = $wb.Worksheets.Item(2)

# Root Extension row: Short/Definition updated from generic placeholder text to the
# specific title/description of this extension.
$ws2.Cells.Item(2, 11).Value = "Medical Fully Insured Indicator"
$ws2.Cells.Item(2, 12).Value = "Indicator of the fully insured medical coverage for the member or employee"
